# ---------------------------------------------------------------------------
# Reproduces the "Add files via upload" re-save of the workshop sign-up sheet:
#   * the header row shrinks from a 16-column (A:P) duplicated grid down to a
#     single 5-column (A:E) grid
#   * the duplicated "Image 1/2/3" + "Text 1/2/3/4" header labels are replaced
#     by three "Option1/Option2/Option3" labels
#   * the thick box border that was drawn around the old header cells is
#     removed
#   * a couple of stray "Yes" marks move/appear in row 3 to line up with the
#     new 5-column layout
#   * workbook-level cosmetic state (selected cell) is refreshed
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the duplicate/overflow columns F:P entirely (nothing left of F
#    shifts) so the grid becomes A:E and the used range/dimension shrinks to
#    A1:E4 automatically.
$ws.Range("F1:P1").EntireColumn.Delete() | Out-Null

# 2) Re-label the three option columns in the header row.
$ws.Range("C1").Value = "Option1"
$ws.Range("D1").Value = "Option2"
$ws.Range("E1").Value = "Option3"

# 3) Remove the thick box border that used to frame the header cells.
$ws.Range("C1:E1").Borders.LineStyle = -4142

# 4) Row 3 now carries "Yes" under Option1 and Option3 (columns C and E),
#    but not Option2 (column D).
$ws.Range("C3").Value = "Yes"
$ws.Range("E3").Value = "Yes"

# 5) Refresh the selected cell shown when the sheet is reopened.
$ws.Range("F10").Select() | Out-Null
